$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) First screenshot drawing (anchorId 3CA0E247): mark the run noProof so it
#    gets <w:rPr><w:noProof/></w:rPr> like the other screenshot runs.
#    (It is the paragraph right before the "... oder nicht." paragraph.)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("oder nicht.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostIdx = $rng.Paragraphs.First.Index
$prevPara = $d.Paragraphs.Item($hostIdx - 1)
$prevPara.Range.NoProofing = 1

# ---------------------------------------------------------------------------
# 2) Insert two new paragraphs right after the "... oder nicht." paragraph:
#    explanatory sentence + a "Source" styled code line "targetCompatibility = 11"
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("oder nicht.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostIdx = $rng.Paragraphs.First.Index
$hostPara = $d.Paragraphs.Item($hostIdx)
$insertAt = $hostPara.Range
$insertAt.Collapse(0)
$insertAt.InsertParagraphAfter() | Out-Null

$newPara1 = $d.Paragraphs.Item($hostIdx + 1)
$xmlFrag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Falls auf dem lokalen System eine höhere Java-Version ("</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>java</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> –</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>version</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">") installiert ist als in der virtuellen Maschine, dann sollte folgende Einstellung ins </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>build.gradle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> übernommen werden.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara1.Range.InsertXML($xmlFrag1) | Out-Null

$newPara2 = $d.Paragraphs.Item($hostIdx + 2)
$xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Source"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>targetCompatibility</w:t></w:r></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 11</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara2.Range.InsertXML($xmlFrag2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Second screenshot drawing (anchorId 2332AF7B): mark the run noProof too.
#    (It is the paragraph right before "Optional auf Kommandozeile:".)
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Optional auf Kommandozeile:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostIdx = $rng.Paragraphs.First.Index
$prevPara = $d.Paragraphs.Item($hostIdx - 1)
$prevPara.Range.NoProofing = 1

# ---------------------------------------------------------------------------
# 4) "gradle task Jar" source line: split into separate spell-checked runs and
#    tag every run (and the paragraph mark) with lang="de-DE"; also tag the
#    following empty Source paragraph mark with lang="de-DE".
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("gradle task Jar", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostIdx = $rng.Paragraphs.First.Index
$taskPara = $d.Paragraphs.Item($hostIdx)
$xmlFrag3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Source"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>gradle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>task</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="de-DE"/></w:rPr><w:t>Jar</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$taskPara.Range.InsertXML($xmlFrag3) | Out-Null

$emptySourcePara = $d.Paragraphs.Item($hostIdx + 1)
$xmlFrag4 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Source"/><w:rPr><w:lang w:val="de-DE"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$emptySourcePara.Range.InsertXML($xmlFrag4) | Out-Null

# ---------------------------------------------------------------------------
# 5) Drop the stray <w:lastRenderedPageBreak/> in front of the "jar" run
#    (sv-SE "jar -xf minimal.jar META-INF/MANIFEST.MF" source line).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("xf minimal.jar META-INF/MANIFEST.MF", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostIdx = $rng.Paragraphs.First.Index
$jarPara = $d.Paragraphs.Item($hostIdx)
$xmlFrag5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Source"/><w:rPr><w:lang w:val="sv-SE"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="000E7B06"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>jar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="000E7B06"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve"> -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="000E7B06"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t>xf</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="000E7B06"><w:rPr><w:lang w:val="sv-SE"/></w:rPr><w:t xml:space="preserve"> minimal.jar META-INF/MANIFEST.MF</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$jarPara.Range.InsertXML($xmlFrag5) | Out-Null
